{"js": "// Remove every paragraph after the title + blank line (the TBV/PUMP tag\n// list and its detail paragraphs), leaving just the \"TBV Tags\" title and\n// the following blank paragraph before the section break.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Keep paragraph 0 (\"TBV Tags\" title) and paragraph 1 (blank line);\n// delete everything from paragraph 2 through the end of the body.\nfor (let i = paragraphs.items.length - 1; i >= 2; i--) {\n  paragraphs.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove every paragraph after the title + blank line (the TBV/PUMP tag\n# list and its detail paragraphs), leaving just the \"TBV Tags\" title and\n# the following blank paragraph before the section break.\n$d = $word.ActiveDocument\n\n# Walk backwards so deleting a paragraph doesn't shift the indices of the\n# ones still to be removed. Paragraph 1 is the title, paragraph 2 is the\n# blank line directly under it - keep both and delete paragraphs 3..Count.\nfor ($i = $d.Paragraphs.Count; $i -ge 3; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
